$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Approved" header column
$ws.Range("K1").Value = "Approved"

# Row 9 - new item, no Approved value set
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "test"
$ws.Range("C9").Value = "test"
$ws.Range("D9").Value = "other"
$ws.Range("E9").Value = "Excellent"
$ws.Range("F9").Value = "Taken"
$ws.Range("G9").Value = "admin.mike@lsu.edu"
$ws.Range("H9").Value = "photo-1763865475041-04mpl.jpg"
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "8"

# Row 10 - new item, Approved = TRUE
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "admin test"
$ws.Range("C10").Value = "This is a test to see if the admin can aprove"
$ws.Range("D10").Value = "other"
$ws.Range("E10").Value = "Poor"
$ws.Range("F10").Value = "Available"
$ws.Range("G10").Value = "admin.mike@lsu.edu"
$ws.Range("H10").Value = "photo-1763868107950-0kdjg.jpg"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "9"
$ws.Range("K10").Value = $true

# Row 11 - new item, Approved = TRUE
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Lamp"
$ws.Range("C11").Value = "A stock image of a lamp"
$ws.Range("D11").Value = "furniture"
$ws.Range("E11").Value = "Excellent"
$ws.Range("F11").Value = "Available"
$ws.Range("G11").Value = "alice@lsu.edu"
$ws.Range("H11").Value = "photo-1763869213903-li17i.jpg"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "10"
$ws.Range("K11").Value = $true
